$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "307.41"
Set-TextValue $ws.Range("E2") "0.16%"
Set-TextValue $ws.Range("E3") "2.85%"
Set-TextValue $ws.Range("D4") "5.142"
Set-TextValue $ws.Range("E4") "2.08%"
Set-TextValue $ws.Range("D5") "0.07610"
Set-TextValue $ws.Range("E5") "-0.82%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D6") "1.628"
Set-TextValue $ws.Range("E6") "0.89%"
$ws.Range("B7").Value = "BTSEToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D7") "2.492"
Set-TextValue $ws.Range("E7") "-1.52%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D8") "0.9003"
Set-TextValue $ws.Range("E8") "1.25%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D9") "0.1111"
Set-TextValue $ws.Range("E9") "10.46%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D10") "0.1773"
Set-TextValue $ws.Range("E10") "2.74%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D11") "0.09188"
Set-TextValue $ws.Range("E11") "3.45%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D12") "0.04183"
Set-TextValue $ws.Range("E12") "-4.63%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D13") "0.1049"
Set-TextValue $ws.Range("E13") "-0.54%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D14") "0.001260"
Set-TextValue $ws.Range("E14") "-0.73%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D15") "0.005831"
Set-TextValue $ws.Range("E15") "0.45%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D16") "3.357"
Set-TextValue $ws.Range("E16") "0.03%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D17") "4.265"
Set-TextValue $ws.Range("E17") "0.77%"
Set-TextValue $ws.Range("E18") "-1.95%"
Set-TextValue $ws.Range("D19") "6.537"
Set-TextValue $ws.Range("E19") "-7.57%"
Set-TextValue $ws.Range("D20") "0.1360"
Set-TextValue $ws.Range("E20") "1.41%"
Set-TextValue $ws.Range("E21") "-16.37%"
Set-TextValue $ws.Range("D22") "0.04136"
Set-TextValue $ws.Range("E22") "-2.02%"
Set-TextValue $ws.Range("E23") "2.31%"
Set-TextValue $ws.Range("D24") "0.004001"
Set-TextValue $ws.Range("E24") "-1.49%"
Set-TextValue $ws.Range("D25") "0.0001300"
Set-TextValue $ws.Range("E25") "6.40%"
Set-TextValue $ws.Range("E38") "2.32%"
Set-TextValue $ws.Range("D39") "0.05184"
Set-TextValue $ws.Range("E39") "0.54%"
Set-TextValue $ws.Range("D40") "0.007768"
Set-TextValue $ws.Range("E40") "-2.27%"
Set-TextValue $ws.Range("D41") "0.1299"
Set-TextValue $ws.Range("E41") "-1.63%"
Set-TextValue $ws.Range("D42") "0.006969"
Set-TextValue $ws.Range("E42") "6.25%"
Set-TextValue $ws.Range("D43") "0.001951"
Set-TextValue $ws.Range("E43") "-1.90%"
Set-TextValue $ws.Range("D44") "0.007617"
Set-TextValue $ws.Range("E44") "-12.14%"
Set-TextValue $ws.Range("D45") "0.3054"
Set-TextValue $ws.Range("E45") "0.35%"
Set-TextValue $ws.Range("D46") "0.00006742"
Set-TextValue $ws.Range("E46") "2.63%"
Set-TextValue $ws.Range("E47") "-0.14%"
Set-TextValue $ws.Range("D48") "0.03134"
Set-TextValue $ws.Range("E48") "820.67%"
Set-TextValue $ws.Range("E50") "-0.14%"
Set-TextValue $ws.Range("E51") "-0.14%"

Write-Host "Applied all changes"